$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---
$ws.Range("D2").Value = "30.762.35"

# --- Row 3 (Ethereum) ---
$ws.Range("D3").Value = "1.885.16"
$ws.Range("E3").Value = "  +2.01%  "

# --- Row 4 (TetherUSD) ---
$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.16%  "

# --- Row 5 (BNB) ---
$ws.Range("D5").Value = "'238.67"
$ws.Range("E5").Value = "  +2.24%  "

# --- Row 6 (USDC) ---
$ws.Range("D6").Value = "'0.9992"
$ws.Range("E6").Value = "  -0.10%  "

# --- Row 7 (XRP) ---
$ws.Range("D7").Value = "'0.4755"
$ws.Range("E7").Value = "  +1.78%  "

# --- Row 8 (Cardano) ---
$ws.Range("D8").Value = "'0.2867"
$ws.Range("E8").Value = "  +5.01%  "

# --- Row 9 (Dogecoin) ---
$ws.Range("D9").Value = "'0.06559"
$ws.Range("E9").Value = "  +4.24%  "

# --- Row 10 (Solana) ---
$ws.Range("D10").Value = "'18.84"
$ws.Range("E10").Value = "  +15.72%  "

# --- Row 11 (WrappedEther) ---
$ws.Range("D11").Value = "1.878.45"
$ws.Range("E11").Value = "  +1.72%  "

# --- Row 12 (Litecoin) ---
$ws.Range("D12").Value = "'96.34"
$ws.Range("E12").Value = "  +14.78%  "

# --- Row 13 (TRON) ---
$ws.Range("D13").Value = "'0.07572"

# --- Row 14 (Polkadot) ---
$ws.Range("D14").Value = "'5.116"
$ws.Range("E14").Value = "  +3.66%  "

# --- Row 15 (Polygon) ---
$ws.Range("D15").Value = "'0.6547"
$ws.Range("E15").Value = "  +5.59%  "

# --- Row 16 (BitcoinCash) ---
$ws.Range("D16").Value = "'308.50"
$ws.Range("E16").Value = "  +34.72%  "

# --- Row 17 (WrappedBTC) ---
$ws.Range("D17").Value = "30.759.16"
$ws.Range("E17").Value = "  +1.35%  "

# --- Row 18 (Avalanche) ---
$ws.Range("E18").Value = "  +6.47%  "

# --- Row 19 (Dai) ---
$ws.Range("D19").Value = "'0.9994"

# --- Row 20 (ShibaInu) ---
$ws.Range("D20").Value = "'0.000007569"
$ws.Range("E20").Value = "  +3.50%  "

# --- Row 21 (WrappedliquidstakedEther2.0) ---
$ws.Range("D21").Value = "2.123.03"
$ws.Range("E21").Value = "  +2.67%  "

# --- Row 22 (BinanceUSD) ---
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.31%  "

# --- Row 23 (Uniswap) ---
$ws.Range("D23").Value = "'5.121"
$ws.Range("E23").Value = "  +4.01%  "

# --- Row 24 (Chainlink) ---
$ws.Range("D24").Value = "'6.166"
$ws.Range("E24").Value = "  +5.06%  "

# --- Row 25 (Cosmos) ---
$ws.Range("D25").Value = "'9.279"
$ws.Range("E25").Value = "  +1.40%  "

# --- Row 26 (Monero) ---
$ws.Range("D26").Value = "'166.50"
$ws.Range("E26").Value = "  +0.43%  "

# --- Row 27 (EthereumClassic) ---
$ws.Range("D27").Value = "'20.13"
$ws.Range("E27").Value = "  +12.82%  "

# --- Row 28 (LidoDAOToken) ---
$ws.Range("D28").Value = "'1.950"
$ws.Range("E28").Value = "  +4.29%  "

# --- Row 29 (Stellar) ---
$ws.Range("D29").Value = "'0.1074"
$ws.Range("E29").Value = "  +5.35%  "

# --- Row 30 (Toncoin) ---
$ws.Range("E30").Value = "  -1.58%  "

# --- Row 31 (InternetComputer(DFINITY)) ---
$ws.Range("D31").Value = "'4.158"
$ws.Range("E31").Value = "  +1.67%  "

# --- Row 32 (Filecoin) ---
$ws.Range("D32").Value = "'3.966"
$ws.Range("E32").Value = "  +4.00%  "

# --- Row 33 (Hedera) ---
$ws.Range("D33").Value = "'0.05044"
$ws.Range("E33").Value = "  +3.43%  "

# --- Row 34 (ARBITRUM) ---
$ws.Range("E34").Value = "  +2.82%  "

# --- Row 35 (ImmutableX) ---
$ws.Range("D35").Value = "'0.7315"
$ws.Range("E35").Value = "  +4.15%  "

# --- Row 36 (HuobiToken) ---
$ws.Range("D36").Value = "'2.710"
$ws.Range("E36").Value = "  +0.83%  "

# --- Row 37 (VeChain) ---
$ws.Range("D37").Value = "'0.01944"
$ws.Range("E37").Value = "  +1.56%  "

# --- Row 38 (MXToken) ---
$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  +1.31%  "

# --- Row 39 (RenderToken) ---
$ws.Range("D39").Value = "'2.074"
$ws.Range("E39").Value = "  +7.35%  "

# --- Row 40 (TrustWalletToken) ---
$ws.Range("D40").Value = "'0.9015"
$ws.Range("E40").Value = "  +4.46%  "

# --- Row 41 (Quant) ---
$ws.Range("D41").Value = "'107.77"
$ws.Range("E41").Value = "  +2.09%  "

# --- Row 43 (TheSandbox) ---
$ws.Range("D43").Value = "'0.4216"
$ws.Range("E43").Value = "  +4.69%  "

# --- Row 44 (FraxShare) ---
$ws.Range("D44").Value = "'5.630"
$ws.Range("E44").Value = "  +2.06%  "

# --- Row 45 (Aave) ---
$ws.Range("D45").Value = "'65.91"
$ws.Range("E45").Value = "  +7.14%  "

# --- Row 46 (Aptos) ---
$ws.Range("D46").Value = "'7.354"
$ws.Range("E46").Value = "  +4.02%  "

# --- Row 47: Algorand -> EnergySwap ---
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'9.030"
$ws.Range("E47").Value = "  +4.45%  "

# --- Row 48: EnergySwap -> Algorand ---
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1225"
$ws.Range("E48").Value = "  +1.84%  "

# --- Row 49 (Elrond) ---
$ws.Range("D49").Value = "'34.78"
$ws.Range("E49").Value = "  +4.33%  "

# --- Row 50 (Cronos) ---
$ws.Range("D50").Value = "'0.05608"
$ws.Range("E50").Value = "  +1.33%  "

# --- Row 51: NEARProtocol -> Decentraland ---
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3847"
$ws.Range("E51").Value = "  +5.40%  "
